# Apply the "Use NatMet high-complexity dataset for ROC calculations" edit.
#
# Summary of the change:
#  - The "mock" sheet gains real smplH (high-complexity native sample)
#    read counts in column J (previously all placeholder zeros), and the
#    D14 value used by the ROC totals is corrected from 4 to 0. The
#    J1/A266 header & footer labels, which are shared-string lookups,
#    follow automatically once the underlying text is updated elsewhere.
#  - The "key" sheet's NATIVES formula text ["smpl_high"] is renamed to
#    ["smplH"] for every data row (23:263).
#  - The "types" sheet's description of the high-complexity native taxa
#    is reworded/re-wrapped across D12:D13.
#  - Selections on all three sheets are moved to reflect where the user
#    was last working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "mock": new smplH counts in column J, and a corrected D14.
# ---------------------------------------------------------------
$mock = $wb.Worksheets.Item("mock")

$mock.Cells.Item(1, 10).Value = "smplH"     # J1 header (was "smpl_high")
$mock.Cells.Item(2, 10).Value = 40000       # J2
$mock.Cells.Item(3, 10).Value = 10000       # J3
$mock.Cells.Item(4, 10).Value = 5000        # J4
$mock.Cells.Item(5, 10).Value = 5000        # J5
$mock.Cells.Item(6, 10).Value = 4           # J6
$mock.Cells.Item(8, 10).Value = 5000        # J8
$mock.Cells.Item(9, 10).Value = 4000        # J9
$mock.Cells.Item(10, 10).Value = 20         # J10
$mock.Cells.Item(11, 10).Value = 103        # J11
$mock.Cells.Item(12, 10).Value = 5000       # J12
$mock.Cells.Item(14, 4).Value = 0           # D14 (was 4)
$mock.Cells.Item(20, 10).Value = 2995       # J20

# ---------------------------------------------------------------
# Sheet "key": rename the NATIVES sample tag for every row 23-263.
# ---------------------------------------------------------------
$key = $wb.Worksheets.Item("key")
$key.Range("E23:E263").Value = "[""smplH""]"

# ---------------------------------------------------------------
# Sheet "types": reword the high-complexity dataset description.
# ---------------------------------------------------------------
$types = $wb.Worksheets.Item("types")
$types.Range("D12").Value = "native taxa from NatMet high-complexity dataset"
$types.Range("D13").Value = "for metagenomics benchmarks (Scyrba et al., 2017)"

# ---------------------------------------------------------------
# Restore the on-screen selections the author left behind, making
# sure the "mock" tab ends up the active/selected one again.
# ---------------------------------------------------------------
$key.Range("E23:E263").Select()
$types.Range("D14").Select()
$mock.Activate()
$mock.Range("D17").Select()
